# Updates the Ligue 1 2023-2024 sheet:
#  - Row 42 and row 43 had their match data (columns F:V) swapped back to the
#    correct fixtures (Reims-Brest now on row 42, Strasbourg-Montpellier on row 43).
#  - Row 51 and row 52 had their match data (columns F:V) swapped back to the
#    correct fixtures (Le Havre-Clermont now on row 51, Lens-Toulouse on row 52).
#  - A new row 90 was appended with the Paris SG - Montpellier result.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 42: Reims vs Brest -------------------------------------------------
$ws.Cells.Item(42, 6).Value = "Reims"
$ws.Cells.Item(42, 7).Value = 1
$ws.Cells.Item(42, 8).Value = "Brest"
$ws.Cells.Item(42, 9).Value = 2
$ws.Cells.Item(42, 10).Value = 1.71
$ws.Cells.Item(42, 11).Value = "28/08/2023 16:01"
$ws.Cells.Item(42, 12).Value = 2.03
$ws.Cells.Item(42, 13).Value = "17/09/2023 14:50"
$ws.Cells.Item(42, 14).Value = 3.98
$ws.Cells.Item(42, 15).Value = "28/08/2023 16:01"
$ws.Cells.Item(42, 16).Value = 3.6
$ws.Cells.Item(42, 17).Value = "17/09/2023 14:53"
$ws.Cells.Item(42, 18).Value = 5.01
$ws.Cells.Item(42, 19).Value = "28/08/2023 16:01"
$ws.Cells.Item(42, 20).Value = 3.96
$ws.Cells.Item(42, 21).Value = "17/09/2023 14:53"
$ws.Cells.Item(42, 22).Value = "https://www.betexplorer.com/football/france/ligue-1/reims-brest/pn1vspJg/"

# --- Row 43: Strasbourg vs Montpellier --------------------------------------
$ws.Cells.Item(43, 6).Value = "Strasbourg"
$ws.Cells.Item(43, 7).Value = 2
$ws.Cells.Item(43, 8).Value = "Montpellier"
$ws.Cells.Item(43, 9).Value = 2
$ws.Cells.Item(43, 10).Value = 2.02
$ws.Cells.Item(43, 11).Value = "28/08/2023 16:01"
$ws.Cells.Item(43, 12).Value = 3.1
$ws.Cells.Item(43, 13).Value = "17/09/2023 14:58"
$ws.Cells.Item(43, 14).Value = 3.61
$ws.Cells.Item(43, 15).Value = "28/08/2023 16:01"
$ws.Cells.Item(43, 16).Value = 3.37
$ws.Cells.Item(43, 17).Value = "17/09/2023 14:31"
$ws.Cells.Item(43, 18).Value = 3.8
$ws.Cells.Item(43, 19).Value = "28/08/2023 16:01"
$ws.Cells.Item(43, 20).Value = 2.44
$ws.Cells.Item(43, 21).Value = "17/09/2023 14:57"
$ws.Cells.Item(43, 22).Value = "https://www.betexplorer.com/football/france/ligue-1/strasbourg-montpellier/fJq2dPIt/"

# --- Row 51: Le Havre vs Clermont -------------------------------------------
$ws.Cells.Item(51, 6).Value = "Le Havre"
$ws.Cells.Item(51, 7).Value = 2
$ws.Cells.Item(51, 8).Value = "Clermont"
$ws.Cells.Item(51, 9).Value = 1
$ws.Cells.Item(51, 10).Value = 2.72
$ws.Cells.Item(51, 11).Value = "11/09/2023 13:39"
$ws.Cells.Item(51, 12).Value = 2.56
$ws.Cells.Item(51, 13).Value = "24/09/2023 14:58"
$ws.Cells.Item(51, 14).Value = 3.25
$ws.Cells.Item(51, 15).Value = "11/09/2023 13:39"
$ws.Cells.Item(51, 16).Value = 3.08
$ws.Cells.Item(51, 17).Value = "24/09/2023 14:58"
$ws.Cells.Item(51, 18).Value = 2.64
$ws.Cells.Item(51, 19).Value = "11/09/2023 13:39"
$ws.Cells.Item(51, 20).Value = 3.26
$ws.Cells.Item(51, 21).Value = "24/09/2023 14:58"
$ws.Cells.Item(51, 22).Value = "https://www.betexplorer.com/football/france/ligue-1/le-havre-clermont/MutixoyI/"

# --- Row 52: Lens vs Toulouse ------------------------------------------------
$ws.Cells.Item(52, 6).Value = "Lens"
$ws.Cells.Item(52, 7).Value = 2
$ws.Cells.Item(52, 8).Value = "Toulouse"
$ws.Cells.Item(52, 9).Value = 1
$ws.Cells.Item(52, 10).Value = 1.58
$ws.Cells.Item(52, 11).Value = "05/09/2023 12:01"
$ws.Cells.Item(52, 12).Value = 1.45
$ws.Cells.Item(52, 13).Value = "24/09/2023 14:59"
$ws.Cells.Item(52, 14).Value = 4.44
$ws.Cells.Item(52, 15).Value = "05/09/2023 12:01"
$ws.Cells.Item(52, 16).Value = 4.73
$ws.Cells.Item(52, 17).Value = "24/09/2023 14:59"
$ws.Cells.Item(52, 18).Value = 5.6
$ws.Cells.Item(52, 19).Value = "05/09/2023 12:01"
$ws.Cells.Item(52, 20).Value = 8.07
$ws.Cells.Item(52, 21).Value = "24/09/2023 14:59"
$ws.Cells.Item(52, 22).Value = "https://www.betexplorer.com/football/france/ligue-1/lens-toulouse/QDnrv765/"

# --- New row 90: Paris SG vs Montpellier ------------------------------------
# Seed the new row's formatting (bold/bordered index cell, datetime-formatted
# date cell, etc.) by copying the row above, then overwrite with new values.
$ws.Range("A89:V89").Copy($ws.Range("A90:V90"))

$ws.Cells.Item(90, 1).Value = 89
$ws.Cells.Item(90, 2).Value = "france"
$ws.Cells.Item(90, 3).Value = "ligue-1"
$ws.Cells.Item(90, 4).Value = "2023-2024"
$ws.Cells.Item(90, 5).Value = 45233.875
$ws.Cells.Item(90, 6).Value = "Paris SG"
$ws.Cells.Item(90, 7).Value = 3
$ws.Cells.Item(90, 8).Value = "Montpellier"
$ws.Cells.Item(90, 9).Value = 0
$ws.Cells.Item(90, 10).Value = 1.32
$ws.Cells.Item(90, 11).Value = "22/10/2023 12:02"
$ws.Cells.Item(90, 12).Value = 1.26
$ws.Cells.Item(90, 13).Value = "03/11/2023 20:55"
$ws.Cells.Item(90, 14).Value = 5.73
$ws.Cells.Item(90, 15).Value = "22/10/2023 12:02"
$ws.Cells.Item(90, 16).Value = 6.79
$ws.Cells.Item(90, 17).Value = "03/11/2023 20:57"
$ws.Cells.Item(90, 18).Value = 7.59
$ws.Cells.Item(90, 19).Value = "22/10/2023 12:02"
$ws.Cells.Item(90, 20).Value = 10.44
$ws.Cells.Item(90, 21).Value = "03/11/2023 20:57"
$ws.Cells.Item(90, 22).Value = "https://www.betexplorer.com/football/france/ligue-1/paris-sg-montpellier/vgYlAVn1/"
